$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2160493827160494
$ws.Range("C2").Value = 0.5154320987654321
$ws.Range("J2").Value = 0.0308641975308642
$ws.Range("P2").Value = 0.1512345679012346
$ws.Range("S2").Value = 0.08641975308641975
$ws.Range("B3").Value = 0.01162790697674419
$ws.Range("C3").Value = 0.02906976744186046
$ws.Range("J3").Value = 0.01744186046511628
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.1918604651162791
$ws.Range("J4").Value = 0.02083333333333333
$ws.Range("P4").Value = 0.625
$ws.Range("S4").Value = 0.3541666666666667
$ws.Range("J5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.05726872246696035
$ws.Range("D6").Value = 0.013215859030837
$ws.Range("F6").Value = 0.07929515418502203
$ws.Range("J6").Value = 0.1938325991189427
$ws.Range("O6").Value = 0.02202643171806168
$ws.Range("Q6").Value = 0.2158590308370044
$ws.Range("R6").Value = 0.06167400881057269
$ws.Range("S6").Value = 0.3568281938325991
$ws.Range("B7").Value = 0.1055900621118012
$ws.Range("D7").Value = 0.02484472049689441
$ws.Range("F7").Value = 0.06211180124223602
$ws.Range("J7").Value = 0.1055900621118012
$ws.Range("O7").Value = 0.03726708074534162
$ws.Range("Q7").Value = 0.1863354037267081
$ws.Range("R7").Value = 0.08074534161490683
$ws.Range("S7").Value = 0.3975155279503105
$ws.Range("B8").Value = 0.08977035490605428
$ws.Range("D8").Value = 0.01670146137787056
$ws.Range("F8").Value = 0.04801670146137787
$ws.Range("J8").Value = 0.1252609603340292
$ws.Range("O8").Value = 0.01670146137787056
$ws.Range("Q8").Value = 0.2025052192066806
$ws.Range("R8").Value = 0.1002087682672234
$ws.Range("S8").Value = 0.4008350730688935
$ws.Range("B9").Value = 0.1327014218009479
$ws.Range("D9").Value = 0.02369668246445497
$ws.Range("E9").Value = 0.004739336492890996
$ws.Range("F9").Value = 0.07582938388625593
$ws.Range("J9").Value = 0.1279620853080569
$ws.Range("O9").Value = 0.004739336492890996
$ws.Range("Q9").Value = 0.1990521327014218
$ws.Range("R9").Value = 0.05687203791469194
$ws.Range("S9").Value = 0.3744075829383886
$ws.Range("B10").Value = 0.10446735395189
$ws.Range("D10").Value = 0.02061855670103093
$ws.Range("E10").Value = 0.0006872852233676976
$ws.Range("F10").Value = 0.06804123711340206
$ws.Range("J10").Value = 0.1175257731958763
$ws.Range("O10").Value = 0.02405498281786942
$ws.Range("Q10").Value = 0.2302405498281787
$ws.Range("R10").Value = 0.0859106529209622
$ws.Range("S10").Value = 0.3484536082474227
$ws.Range("F11").Value = 0.003937007874015748
$ws.Range("G11").Value = 0.1417322834645669
$ws.Range("J11").Value = 0.07874015748031496
$ws.Range("K11").Value = 0.1771653543307087
$ws.Range("L11").Value = 0.5748031496062992
$ws.Range("S11").Value = 0.02362204724409449
$ws.Range("G12").Value = 0.7241379310344828
$ws.Range("J12").Value = 0.2206896551724138
$ws.Range("K12").Value = 0.01379310344827586
$ws.Range("L12").Value = 0.02758620689655172
$ws.Range("S12").Value = 0.01379310344827586
$ws.Range("G13").Value = 0.6571428571428571
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.05714285714285714
$ws.Range("F15").Value = 0.01945525291828794
$ws.Range("H15").Value = 0.1673151750972763
$ws.Range("I15").Value = 0.03501945525291829
$ws.Range("J15").Value = 0.4007782101167315
$ws.Range("K15").Value = 0.0311284046692607
$ws.Range("M15").Value = 0.007782101167315175
$ws.Range("O15").Value = 0.05836575875486381
$ws.Range("S15").Value = 0.2801556420233463
$ws.Range("F16").Value = 0.02912621359223301
$ws.Range("H16").Value = 0.1844660194174757
$ws.Range("I16").Value = 0.05339805825242718
$ws.Range("J16").Value = 0.441747572815534
$ws.Range("K16").Value = 0.06796116504854369
$ws.Range("M16").Value = 0.01941747572815534
$ws.Range("O16").Value = 0.07766990291262135
$ws.Range("S16").Value = 0.1262135922330097
$ws.Range("F17").Value = 0.009107468123861567
$ws.Range("H17").Value = 0.1730418943533698
$ws.Range("I17").Value = 0.1056466302367942
$ws.Range("J17").Value = 0.4389799635701275
$ws.Range("K17").Value = 0.09107468123861566
$ws.Range("M17").Value = 0.02185792349726776
$ws.Range("O17").Value = 0.0692167577413479
$ws.Range("S17").Value = 0.09107468123861566
$ws.Range("F18").Value = 0.004739336492890996
$ws.Range("H18").Value = 0.1563981042654028
$ws.Range("I18").Value = 0.06635071090047394
$ws.Range("J18").Value = 0.5213270142180095
$ws.Range("K18").Value = 0.06635071090047394
$ws.Range("M18").Value = 0.02369668246445497
$ws.Range("N18").Value = 0.004739336492890996
$ws.Range("O18").Value = 0.07582938388625593
$ws.Range("S18").Value = 0.08056872037914692
$ws.Range("F19").Value = 0.00906344410876133
$ws.Range("H19").Value = 0.2046827794561933
$ws.Range("I19").Value = 0.08987915407854985
$ws.Range("J19").Value = 0.4025679758308157
$ws.Range("K19").Value = 0.09063444108761329
$ws.Range("M19").Value = 0.01057401812688822
$ws.Range("N19").Value = 0.001510574018126888
$ws.Range("O19").Value = 0.06570996978851963
$ws.Range("S19").Value = 0.1253776435045317
